# "Generate Report for Handoff"
#
# The localization-status report is being refreshed: the overall
# status moves from "In Translation" to "Ready for handoff", and the
# handoff timestamps are bumped forward to the moment the new Xliff
# files were produced. Updating the text also makes the Status /
# Latest Handoff Datetime columns a bit wider so the new values aren't
# clipped (mirrors what Excel's own column auto-fit would do).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------
# This shared string is used by the Status-ish cells on all three
# sheets; update every cell that shows it so the workbook ends up with
# a single, de-duplicated "Ready for handoff" string.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Handoff datetimes ------------------------------------------------
# zh-cn's handoff xliff was (re)generated at 18:47:49 ...
$wsZhCn.Range("H2").Value = "2016-08-19 18:47:49"

# ... while de-de's handoff (surfaced both on the de-de sheet and
# summarized on the Overview sheet) landed a few seconds later, at
# 18:47:53.
$wsDeDe.Range("H2").Value = "2016-08-19 18:47:53"
$wsOverview.Range("G2").Value = "2016-08-19 18:47:53"

# --- Column widths ------------------------------------------------------
# Widen the Status columns (Overview E:F) and the Status column on each
# language sheet (column C) to fit the longer "Ready for handoff" text.
$wsOverview.Range("E1:F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333332
